$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.041845560073853
$ws.Range("B1").Value = 3.458868503570557
$ws.Range("C1").Value = 3.169588804244995
$ws.Range("D1").Value = 3.60242772102356
$ws.Range("E1").Value = 1.498715281486511
